$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 50006500
$ws.Range("J32").Value = 7220.6665
$ws.Range("L32").Value = 7220.6665
$ws.Range("N32").Value = -7872.6665

$ws.Range("H40").Value = 90551.836
$ws.Range("J40").Value = 88662.2
$ws.Range("L40").Value = 88662.2
$ws.Range("N40").Value = -89012.2

$ws.Range("H64").Value = 3889
$ws.Range("J64").Value = 3912.5
$ws.Range("L64").Value = 3912.5
$ws.Range("N64").Value = -4408.5

$ws.Range("H67").Value = 3889
$ws.Range("J67").Value = 3912.5
$ws.Range("L67").Value = 3912.5
$ws.Range("N67").Value = -5628.5

$ws.Range("H70").Value = 8623.154
$ws.Range("J70").Value = 9033.223
$ws.Range("L70").Value = 27099.669
$ws.Range("N70").Value = -27639.669

$ws.Range("H73").Value = 8623.154
$ws.Range("J73").Value = 9033.223
$ws.Range("L73").Value = 27099.669
$ws.Range("N73").Value = -28971.669

$ws.Range("H76").Value = 1724.5
$ws.Range("I76").Value = 450
$ws.Range("J76").Value = 2999
$ws.Range("K76").Value = 450
$ws.Range("L76").Value = 2999
$ws.Range("M76").Value = -135
$ws.Range("N76").Value = -3629

$ws.Range("H79").Value = 1724.5
$ws.Range("I79").Value = 450
$ws.Range("J79").Value = 2999
$ws.Range("K79").Value = 450
$ws.Range("L79").Value = 2999
$ws.Range("M79").Value = 642
$ws.Range("N79").Value = -5183

$ws.Range("H86").Value = 250003500
$ws.Range("I86").Value = 500002500
$ws.Range("J86").Value = 4499.5
$ws.Range("K86").Value = 500002500
$ws.Range("L86").Value = 4499.5
$ws.Range("M86").Value = -500001377
$ws.Range("N86").Value = -6745.5

$ws.Range("H89").Value = 250003500
$ws.Range("I89").Value = 500002500
$ws.Range("J89").Value = 4499.5
$ws.Range("K89").Value = 2500012500
$ws.Range("L89").Value = 22497.5
$ws.Range("M89").Value = -2500006884
$ws.Range("N89").Value = -33729.5

$ws.Range("H132").Value = 10317.762
$ws.Range("I132").Value = 2451.9412
$ws.Range("K132").Value = 7355.823600000001
$ws.Range("M132").Value = -4825.823600000001

$ws.Range("H137").Value = 2918.4119
$ws.Range("I137").Value = 2701
$ws.Range("K137").Value = 8103
$ws.Range("M137").Value = -5553

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 968.7143
$ws.Range("I4").Value = 963.6667
$ws.Range("K4").Value = 963.6667
$ws.Range("M4").Value = -847.6667

$ws.Range("H32").Value = 10951.326
$ws.Range("I32").Value = 10054.53
$ws.Range("K32").Value = 10054.53
$ws.Range("M32").Value = -9767.530000000001

$ws.Range("H45").Value = 2925
$ws.Range("I45").Value = 1960
$ws.Range("J45").Value = 3614.2856
$ws.Range("K45").Value = 1960
$ws.Range("L45").Value = 3614.2856
$ws.Range("M45").Value = -1583
$ws.Range("N45").Value = -4368.2856

$ws.Range("H61").Value = 5937.75
$ws.Range("I61").Value = 4932.8335
$ws.Range("K61").Value = 4932.8335
$ws.Range("M61").Value = -4720.8335

$ws.Range("H74").Value = 1286.4286
$ws.Range("I74").Value = 1217.4615
$ws.Range("K74").Value = 1217.4615
$ws.Range("M74").Value = -343.4614999999999

$ws.Range("H77").Value = 1286.4286
$ws.Range("I77").Value = 1217.4615
$ws.Range("K77").Value = 6087.307499999999
$ws.Range("M77").Value = -1719.307499999999

$ws.Range("H110").Value = 2451.5
$ws.Range("I110").Value = 2155.5
$ws.Range("J110").Value = 2747.5
$ws.Range("K110").Value = 2155.5
$ws.Range("L110").Value = 2747.5
$ws.Range("M110").Value = -110.5
$ws.Range("N110").Value = -6837.5

$ws.Range("H124").Value = 57500
$ws.Range("J124").Value = 57500
$ws.Range("L124").Value = 57500
$ws.Range("N124").Value = -67320

$ws.Range("H132").Value = 22044.152
$ws.Range("I132").Value = 29614.65
$ws.Range("J132").Value = 10397.23
$ws.Range("K132").Value = 88843.95000000001
$ws.Range("L132").Value = 31191.69
$ws.Range("M132").Value = -86313.95000000001
$ws.Range("N132").Value = -36251.69

$ws.Range("H136").Value = 5937.75
$ws.Range("I136").Value = 4932.8335
$ws.Range("K136").Value = 14798.5005
$ws.Range("M136").Value = -12248.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2826.3333
$ws.Range("J20").Value = 2768.3333
$ws.Range("L20").Value = 2768.3333
$ws.Range("N20").Value = -3262.3333

$ws.Range("H99").Value = 2136.6667
$ws.Range("I99").Value = 2136.6667
$ws.Range("K99").Value = 2136.6667
$ws.Range("M99").Value = -638.6667000000002

$ws.Range("H132").Value = 91709.39999999999
$ws.Range("J132").Value = 91709.39999999999
$ws.Range("L132").Value = 91709.39999999999
$ws.Range("N132").Value = -101829.4

$ws.Range("H134").Value = 4125.6
$ws.Range("I134").Value = 3572.5
$ws.Range("K134").Value = 10717.5
$ws.Range("M134").Value = -8182.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 63511924
$ws.Range("I132").Value = 88895190
$ws.Range("K132").Value = 266685570
$ws.Range("M132").Value = -266683040

$ws.Range("H134").Value = 2975.037
$ws.Range("I134").Value = 2538.7144
$ws.Range("K134").Value = 7616.1432
$ws.Range("M134").Value = -5081.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 69511.14999999999
$ws.Range("K4").Value = 208533.45
$ws.Range("M4").Value = -208421.45

$ws.Range("H92").Value = 1409.6
$ws.Range("I92").Value = 1849.5
$ws.Range("J92").Value = 1116.3334
$ws.Range("K92").Value = 5548.5
$ws.Range("L92").Value = 3349.0002
$ws.Range("M92").Value = -4300.5
$ws.Range("N92").Value = -5845.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6965.4136
$ws.Range("I70").Value = 6789.263
$ws.Range("J70").Value = 7300.1
$ws.Range("K70").Value = 6789.263
$ws.Range("L70").Value = 7300.1
$ws.Range("M70").Value = -6519.263
$ws.Range("N70").Value = -7840.1

$ws.Range("H73").Value = 6965.4136
$ws.Range("I73").Value = 6789.263
$ws.Range("J73").Value = 7300.1
$ws.Range("K73").Value = 6789.263
$ws.Range("L73").Value = 7300.1
$ws.Range("M73").Value = -5853.263
$ws.Range("N73").Value = -9172.1

$ws.Range("H132").Value = 3946.818
$ws.Range("I132").Value = 3770.2666
$ws.Range("J132").Value = 4325.143
$ws.Range("K132").Value = 11310.7998
$ws.Range("L132").Value = 12975.429
$ws.Range("M132").Value = -8780.799800000001
$ws.Range("N132").Value = -18035.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3781.5334
$ws.Range("I7").Value = 3198.6274
$ws.Range("J7").Value = 5020.2085
$ws.Range("K7").Value = 3198.6274
$ws.Range("L7").Value = 5020.2085
$ws.Range("M7").Value = -3086.6274
$ws.Range("N7").Value = -5244.2085

$ws.Range("H16").Value = 3722
$ws.Range("I16").Value = 1699.5
$ws.Range("K16").Value = 1699.5
$ws.Range("M16").Value = -1529.5

$ws.Range("H68").Value = 3887
$ws.Range("J68").Value = 5142.727
$ws.Range("L68").Value = 5142.727
$ws.Range("N68").Value = -6640.727

$ws.Range("H71").Value = 3887
$ws.Range("J71").Value = 5142.727
$ws.Range("L71").Value = 25713.635
$ws.Range("N71").Value = -33201.63499999999

$ws.Range("H82").Value = 1519
$ws.Range("J82").Value = 1348.6
$ws.Range("L82").Value = 1348.6
$ws.Range("N82").Value = -2070.6

$ws.Range("H85").Value = 1519
$ws.Range("J85").Value = 1348.6
$ws.Range("L85").Value = 1348.6
$ws.Range("N85").Value = -3844.6

$ws.Range("H93").Value = 5000
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").Value = ""

$ws.Range("H126").Value = 3781.5334
$ws.Range("I126").Value = 3198.6274
$ws.Range("J126").Value = 5020.2085
$ws.Range("K126").Value = 9595.8822
$ws.Range("L126").Value = 15060.6255
$ws.Range("M126").Value = -7125.8822
$ws.Range("N126").Value = -20000.6255

$ws.Range("H132").Value = 3490.7097
$ws.Range("I132").Value = 2662.7974
$ws.Range("J132").Value = 8162.5
$ws.Range("K132").Value = 7988.3922
$ws.Range("L132").Value = 24487.5
$ws.Range("M132").Value = -5458.3922
$ws.Range("N132").Value = -29547.5

$ws.Range("H136").Value = 4228.436
$ws.Range("I136").Value = 3100.3103
$ws.Range("K136").Value = 9300.930899999999
$ws.Range("M136").Value = -6750.930899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11999.833
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376

$ws.Range("H65").Value = 11999.833
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880

$ws.Range("H132").Value = 1739.4073
$ws.Range("I132").Value = 1551.8823
$ws.Range("J132").Value = 2058.2
$ws.Range("K132").Value = 4655.6469
$ws.Range("L132").Value = 6174.599999999999
$ws.Range("M132").Value = -2125.6469
$ws.Range("N132").Value = -11234.6
